$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Akses" values in column A (A2:A7) - renumber from 7-12 to 1-6
$ws.Range("A2").Value = "Akses1"
$ws.Range("A3").Value = "Akses2"
$ws.Range("A4").Value = "Akses3"
$ws.Range("A5").Value = "Akses4"
$ws.Range("A6").Value = "Akses5"
$ws.Range("A7").Value = "Akses6"

# Update the selection to a single cell C6
$ws.Range("C6").Select()
